$wb = $excel.ActiveWorkbook

# New field-sampling data collected on 2020-09-08 for each logger location.
# Each sheet gets one new row appended (datetime, conductivity_uscm, temp_C),
# then the sheet's selection is moved to mirror the saved worksheet view.

$ws = $wb.Worksheets.Item("WIC")
$ws.Range("A8").Value = 44082.375
$ws.Range("B8").Value = 638.70000000000005
$ws.Range("C8").Value = 17.8
[void]$ws.Range("A8").Select()

$ws = $wb.Worksheets.Item("YS")
$ws.Range("A23").Value = 44082.489583333336
$ws.Range("B23").Value = 58.99
$ws.Range("C23").Value = 17.2
[void]$ws.Range("B29").Select()

$ws = $wb.Worksheets.Item("SW")
$ws.Range("A22").Value = 44082.511805555558
$ws.Range("B22").Value = 911.7
$ws.Range("C22").Value = 16
[void]$ws.Range("A22").Select()

$ws = $wb.Worksheets.Item("YI")
$ws.Range("A21").Value = 44082.39166666667
$ws.Range("B21").Value = 433.4
$ws.Range("C21").Value = 18.8
[void]$ws.Range("C22").Select()

$ws = $wb.Worksheets.Item("YN")
$ws.Range("A21").Value = 44082.409722222219
$ws.Range("B21").Value = 557.4
$ws.Range("C21").Value = 17.2
[void]$ws.Range("A21").Select()

$ws = $wb.Worksheets.Item("6MC")
$ws.Range("A22").Value = 44082.429166666669
$ws.Range("B22").Value = 660.2
$ws.Range("C22").Value = 14.4
[void]$ws.Range("G26").Select()

$ws = $wb.Worksheets.Item("DC")
$ws.Range("A22").Value = 44082.439583333333
$ws.Range("B22").Value = 657.8
$ws.Range("C22").Value = 13.1
[void]$ws.Range("A22").Select()

$ws = $wb.Worksheets.Item("PBMS")
$ws.Range("A23").Value = 44082.45416666667
$ws.Range("B23").Value = 742.7
$ws.Range("C23").Value = 16.100000000000001
[void]$ws.Range("B26").Select()

$ws = $wb.Worksheets.Item("PBSF")
$ws.Range("A24").Value = 44082.463194444441
$ws.Range("B24").Value = 288.3
$ws.Range("C24").Value = 16.7
[void]$ws.Range("A24").Select()
